$d = $word.ActiveDocument

# Locate the paragraph ending in "... FreiePlaetzeAnzeigen erweitert"
# (Jonathan's last bullet). We search by content instead of a hard-coded
# index so the script is resilient to that.
$targetIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $t = $d.Paragraphs.Item($i).Range.Text
    if ($t -like "*FreiePlaetzeAnzeigen erweitert*") {
        $targetIndex = $i
        break
    }
}

$p8 = $d.Paragraphs.Item($targetIndex)

# Create a new paragraph right after it; it automatically inherits that
# paragraph's list/paragraph formatting (pStyle "Listenabsatz", numId 4).
$p8.Range.InsertParagraphAfter()

$newIndex = $targetIndex + 1
$newPara = $d.Paragraphs.Item($newIndex)
$newPara.Range.InsertBefore("Vier verschiedene Methoden zum Filtern von Vorstellungen statt nur einer")

# Re-fetch the paragraph now that it contains text, and compute the
# position right after the inserted text (i.e. right before the
# paragraph mark) -- that's where the _GoBack bookmark must live, since
# it always marks the location of the most recent edit.
$newPara2 = $d.Paragraphs.Item($newIndex)
$textEndPos = $newPara2.Range.End - 1

# Placing a zero-width bookmark exactly at a paragraph-mark position can
# misbehave in this runtime, so as a workaround we temporarily append a
# dummy character right after the text, anchor the bookmark just before
# that dummy character (a safe, non-paragraph-mark position), and then
# remove the dummy character again -- leaving the bookmark correctly
# collapsed right after the real text, matching real Word's behaviour.
$lastCharRange = $d.Range($textEndPos - 1, $textEndPos)
$lastCharRange.InsertAfter("X")

$d.Bookmarks.Item("_GoBack").Delete()
$bookmarkTarget = $d.Range($textEndPos, $textEndPos)
$d.Bookmarks.Add("_GoBack", $bookmarkTarget)

$dummyRange = $d.Range($textEndPos, $textEndPos + 1)
$dummyRange.Delete()
